# "Added last minute updates"
#
# 1. First paragraph: add a paragraph border (top/left/bottom/right) whose
#    only attribute is a 5pt "space" (distance from text) -- no line, just
#    spacing -- and bump the left indent from 120 twips (6pt) to 225 twips
#    (11.25pt).
# 2. Replace the merge-field-style bookmark text in that paragraph's run,
#    and drop the now-unneeded trailing space run that followed it.

$d = $word.ActiveDocument
$p1 = $d.Paragraphs.Item(1)

# --- paragraph border (w:pBdr) with 5pt spacing on all four sides ---
$borders = $p1.Format.Borders
$borders.DistanceFromTop    = 5
$borders.DistanceFromLeft   = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight  = 5

# --- indent: 120 twips -> 225 twips (LeftIndent is in points) ---
$p1.Format.LeftIndent = 11.25

# --- drop the trailing " " run at the end of paragraph 1 (positions are
#     still valid pre-replace, since the replace below changes text length) ---
$spaceRange = $d.Range($p1.Range.End - 2, $p1.Range.End - 1)
if ($spaceRange.Text -eq " ") {
    $spaceRange.Delete()
}

# --- update the bookmark-style placeholder text ---
$d.Content.Find.Execute(
    "**ID__AFFARS_pgi_5342_topic_2__ID**", $true, $false, $false, $false,
    $false, $true, 1, $false, "**ID__AFFARS_AF_PGI_5342__ID**", 2) | Out-Null
